# Adds daily reading report rows for 17-10-10 through 17-10-27 (rows 18-35),
# matching the data already present for earlier dates in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[18] = @{
    A = '17-10-10 ~ 17-10-11'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '5'
    J = '77'
    K = '93.90%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000031'', ''7ff9010202000057'', ''7ff9010202000060'', ''7ff9010202000104'', ''7ff9010202000126'']'
}

$rows[19] = @{
    A = '17-10-11 ~ 17-10-12'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '2'
    J = '80'
    K = '97.56%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000151'', ''7ff9010202000163'']'
}

$rows[20] = @{
    A = '17-10-12 ~ 17-10-13'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '5'
    J = '77'
    K = '93.90%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000059'', ''7ff9010202000082'', ''7ff9010202000097'', ''7ff9010202000106'', ''7ff9010202000123'']'
}

$rows[21] = @{
    A = '17-10-13 ~ 17-10-14'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '5'
    J = '77'
    K = '93.90%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000067'', ''7ff9010202000082'', ''7ff9010202000088'', ''7ff9010202000129'', ''7ff9010202000163'']'
}

$rows[22] = @{
    A = '17-10-14 ~ 17-10-15'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '3'
    J = '79'
    K = '96.34%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000057'', ''7ff9010202000121'', ''7ff9010202000146'']'
}

$rows[23] = @{
    A = '17-10-15 ~ 17-10-16'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[24] = @{
    A = '17-10-16 ~ 17-10-17'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[25] = @{
    A = '17-10-17 ~ 17-10-18'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[26] = @{
    A = '17-10-18 ~ 17-10-19'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[27] = @{
    A = '17-10-19 ~ 17-10-20'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[28] = @{
    A = '17-10-20 ~ 17-10-21'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[29] = @{
    A = '17-10-21 ~ 17-10-22'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[30] = @{
    A = '17-10-22 ~ 17-10-23'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[31] = @{
    A = '17-10-23 ~ 17-10-24'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[32] = @{
    A = '17-10-24 ~ 17-10-25'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[33] = @{
    A = '17-10-25 ~ 17-10-26'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[34] = @{
    A = '17-10-26 ~ 17-10-27'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

$rows[35] = @{
    A = '17-10-27 ~ 17-10-28'
    B = '82'
    C = '82'
    D = '100.00%'
    E = '0'
    F = '82'
    G = '100.00%'
    H = '0'
    I = '82'
    J = '0'
    K = '0.00%'
    L = '[]'
    M = '[]'
    N = '[''7ff9010202000016'', ''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
}

foreach ($r in 18..35) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $row["A"]
    $ws.Cells.Item($r, 2).Value = $row["B"]
    $ws.Cells.Item($r, 3).Value = $row["C"]
    $ws.Cells.Item($r, 4).Value = $row["D"]
    $ws.Cells.Item($r, 5).Value = $row["E"]
    $ws.Cells.Item($r, 6).Value = $row["F"]
    $ws.Cells.Item($r, 7).Value = $row["G"]
    $ws.Cells.Item($r, 8).Value = $row["H"]
    $ws.Cells.Item($r, 9).Value = $row["I"]
    $ws.Cells.Item($r, 10).Value = $row["J"]
    $ws.Cells.Item($r, 11).Value = $row["K"]
    $ws.Cells.Item($r, 12).Value = $row["L"]
    $ws.Cells.Item($r, 13).Value = $row["M"]
    $ws.Cells.Item($r, 14).Value = $row["N"]
}

$ws.Range("A25").Select()
